# Applies the "Penalty Reward System" forecast-shift edit:
#  - Sheet "Forecast Comparison": Week_Start_Date values shift up one row
#    (each row gets the following week's old date, row 17 gets a new
#    extrapolated date), and the MyForecast (D) column is recalculated.
#  - Sheet "Summary": several derived metrics are refreshed to match the
#    new forecast numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Forecast Comparison
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Keep these as plain text (they were stored as text dates, not Excel
# date serials) by forcing a text number format before assigning.
$ws1.Range("B2:B17").NumberFormat = "@"

$weekDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecast = @(16, 16, 16, 16, 16, 16, 16, 17, 17, 17, 17, 17, 15, 15, 15, 15)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $weekDates[$i]
    $ws1.Cells.Item($row, 4).Value = $myForecast[$i]
}

# ---------------------------------------------------------------------
# Sheet 2: Summary
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

# These cells are stored as text (even the numeric-looking ones), so
# force a text format before writing the new values.
$ws2.Range("B2:B15").NumberFormat = "@"

$ws2.Range("B2").Value = "2023-01-22 to 2025-01-05"
$ws2.Range("B4").Value = "44"
$ws2.Range("B8").Value = "1435 units"
$ws2.Range("B9").Value = "258"
$ws2.Range("B10").Value = "130"
$ws2.Range("B11").Value = "64"
$ws2.Range("B12").Value = "17"
$ws2.Range("B13").Value = "2025-03-09"
$ws2.Range("B14").Value = "15"
$ws2.Range("B15").Value = "2025-04-20"
